# Scheduled runner update: refresh market-price / profit figures across all
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with latest Universalis
# market-board snapshot values.

$wb = $excel.ActiveWorkbook

# Sheet 1, Row 64
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(64, 8).Value = 4515.2104
$ws.Cells.Item(64, 9).Value = 3772.5
$ws.Cells.Item(64, 10).Value = 4713.2666
$ws.Cells.Item(64, 11).Value = 3772.5
$ws.Cells.Item(64, 12).Value = 4713.2666
$ws.Cells.Item(64, 13).Value = -3524.5
$ws.Cells.Item(64, 14).Value = -5209.2666

# Sheet 1, Row 67
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(67, 8).Value = 4515.2104
$ws.Cells.Item(67, 9).Value = 3772.5
$ws.Cells.Item(67, 10).Value = 4713.2666
$ws.Cells.Item(67, 11).Value = 3772.5
$ws.Cells.Item(67, 12).Value = 4713.2666
$ws.Cells.Item(67, 13).Value = -2914.5
$ws.Cells.Item(67, 14).Value = -6429.2666

# Sheet 1, Row 82
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(82, 8).Value = 23670.8
$ws.Cells.Item(82, 9).Value = 12076.75
$ws.Cells.Item(82, 11).Value = 36230.25
$ws.Cells.Item(82, 13).Value = -35824.25

# Sheet 1, Row 85
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(85, 8).Value = 23670.8
$ws.Cells.Item(85, 9).Value = 12076.75
$ws.Cells.Item(85, 11).Value = 36230.25
$ws.Cells.Item(85, 13).Value = -34826.25

# Sheet 1, Row 98
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(98, 8).Value = 71431576
$ws.Cells.Item(98, 9).Value = 71431576
$ws.Cells.Item(98, 11).Value = 71431576
$ws.Cells.Item(98, 13).Value = -71430078

# Sheet 1, Row 100
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(100, 8).Value = 3110.389
$ws.Cells.Item(100, 10).Value = 3374.1875
$ws.Cells.Item(100, 12).Value = 3374.1875
$ws.Cells.Item(100, 14).Value = -4456.1875

# Sheet 1, Row 122
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(122, 8).Value = 71431576
$ws.Cells.Item(122, 9).Value = 71431576
$ws.Cells.Item(122, 11).Value = 214294728
$ws.Cells.Item(122, 13).Value = -214292278

# Sheet 1, Row 132
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(132, 8).Value = 3138.6667
$ws.Cells.Item(132, 9).Value = 3142.7715
$ws.Cells.Item(132, 11).Value = 9428.3145
$ws.Cells.Item(132, 13).Value = -6898.3145

# Sheet 1, Row 137
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(137, 8).Value = 3500.6924
$ws.Cells.Item(137, 9).Value = 2320.7908
$ws.Cells.Item(137, 10).Value = 5806.864
$ws.Cells.Item(137, 11).Value = 6962.3724
$ws.Cells.Item(137, 12).Value = 17420.592
$ws.Cells.Item(137, 13).Value = -4412.3724
$ws.Cells.Item(137, 14).Value = -22520.592

# Sheet 1, Row 138
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(138, 8).Value = 2303.5293
$ws.Cells.Item(138, 10).Value = 3026.2341
$ws.Cells.Item(138, 12).Value = 9078.7023
$ws.Cells.Item(138, 14).Value = -19358.7023

# Sheet 1, Row 141
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(141, 8).Value = 3151.3667
$ws.Cells.Item(141, 9).Value = 3385.6296
$ws.Cells.Item(141, 10).Value = 1043
$ws.Cells.Item(141, 11).Value = 10156.8888
$ws.Cells.Item(141, 12).Value = 3129
$ws.Cells.Item(141, 13).Value = -4976.888800000001
$ws.Cells.Item(141, 14).Value = -13489

# Sheet 2, Row 4
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 8).Value = 592
$ws.Cells.Item(4, 9).Value = 775.7143
$ws.Cells.Item(4, 10).Value = 163.33333
$ws.Cells.Item(4, 11).Value = 775.7143
$ws.Cells.Item(4, 12).Value = 163.33333
$ws.Cells.Item(4, 13).Value = -659.7143
$ws.Cells.Item(4, 14).Value = -395.33333

# Sheet 2, Row 28
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(28, 8).Value = 31487
$ws.Cells.Item(28, 9).Value = 11608.75
$ws.Cells.Item(28, 10).Value = 111000
$ws.Cells.Item(28, 11).Value = 11608.75
$ws.Cells.Item(28, 12).Value = 111000
$ws.Cells.Item(28, 13).Value = -11416.75
$ws.Cells.Item(28, 14).Value = -111384

# Sheet 2, Row 45
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(45, 8).Value = 33334810
$ws.Cells.Item(45, 9).Value = 35715724
$ws.Cells.Item(45, 11).Value = 35715724
$ws.Cells.Item(45, 13).Value = -35715347

# Sheet 2, Row 61
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 21434642
$ws.Cells.Item(61, 9).Value = 16133808
$ws.Cells.Item(61, 11).Value = 16133808
$ws.Cells.Item(61, 13).Value = -16133596

# Sheet 2, Row 74
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74, 8).Value = 4818840.5
$ws.Cells.Item(74, 9).Value = 5816046
$ws.Cells.Item(74, 10).Value = 920671.8
$ws.Cells.Item(74, 11).Value = 5816046
$ws.Cells.Item(74, 12).Value = 920671.8
$ws.Cells.Item(74, 13).Value = -5815172
$ws.Cells.Item(74, 14).Value = -922419.8

# Sheet 2, Row 77
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(77, 8).Value = 4818840.5
$ws.Cells.Item(77, 9).Value = 5816046
$ws.Cells.Item(77, 10).Value = 920671.8
$ws.Cells.Item(77, 11).Value = 29080230
$ws.Cells.Item(77, 12).Value = 4603359
$ws.Cells.Item(77, 13).Value = -29075862
$ws.Cells.Item(77, 14).Value = -4612095

# Sheet 2, Row 97
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(97, 8).Value = 1173.4878
$ws.Cells.Item(97, 9).Value = 679.86957
$ws.Cells.Item(97, 11).Value = 679.86957
$ws.Cells.Item(97, 13).Value = -183.86957

# Sheet 2, Row 99
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(99, 8).Value = 31487
$ws.Cells.Item(99, 9).Value = 11608.75
$ws.Cells.Item(99, 10).Value = 111000
$ws.Cells.Item(99, 11).Value = 11608.75
$ws.Cells.Item(99, 12).Value = 111000
$ws.Cells.Item(99, 13).Value = -8613.75
$ws.Cells.Item(99, 14).Value = -116990

# Sheet 2, Row 122
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(122, 8).Value = 4382.4
$ws.Cells.Item(122, 9).Value = 3970.6667
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 11912.0001
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -9462.000100000001
$ws.Cells.Item(122, 14).Value = -19900

# Sheet 2, Row 132
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(132, 8).Value = 4707.5
$ws.Cells.Item(132, 9).Value = 2362
$ws.Cells.Item(132, 10).Value = 8393.286
$ws.Cells.Item(132, 11).Value = 7086
$ws.Cells.Item(132, 12).Value = 25179.858
$ws.Cells.Item(132, 13).Value = -4556
$ws.Cells.Item(132, 14).Value = -30239.858

# Sheet 2, Row 136
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136, 8).Value = 21434642
$ws.Cells.Item(136, 9).Value = 16133808
$ws.Cells.Item(136, 11).Value = 48401424
$ws.Cells.Item(136, 13).Value = -48398874

# Sheet 3, Row 20
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 6563.32
$ws.Cells.Item(20, 9).Value = 6875.375
$ws.Cells.Item(20, 10).Value = 6008.5557
$ws.Cells.Item(20, 11).Value = 6875.375
$ws.Cells.Item(20, 12).Value = 6008.5557
$ws.Cells.Item(20, 13).Value = -6628.375
$ws.Cells.Item(20, 14).Value = -6502.5557

# Sheet 3, Row 86
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(86, 8).Value = 3949.7036
$ws.Cells.Item(86, 9).Value = 4291.45
$ws.Cells.Item(86, 11).Value = 4291.45
$ws.Cells.Item(86, 13).Value = -3168.45

# Sheet 3, Row 89
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(89, 8).Value = 3949.7036
$ws.Cells.Item(89, 9).Value = 4291.45
$ws.Cells.Item(89, 11).Value = 21457.25
$ws.Cells.Item(89, 13).Value = -15841.25

# Sheet 3, Row 99
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(99, 8).Value = 8196.19
$ws.Cells.Item(99, 9).Value = 13491.111
$ws.Cells.Item(99, 10).Value = 4225
$ws.Cells.Item(99, 11).Value = 13491.111
$ws.Cells.Item(99, 12).Value = 4225
$ws.Cells.Item(99, 13).Value = -11993.111
$ws.Cells.Item(99, 14).Value = -7221

# Sheet 3, Row 107
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(107, 8).Value = 617.087
$ws.Cells.Item(107, 9).Value = 578.8421
$ws.Cells.Item(107, 10).Value = 798.75
$ws.Cells.Item(107, 11).Value = 578.8421
$ws.Cells.Item(107, 12).Value = 798.75
$ws.Cells.Item(107, 13).Value = 1341.1579
$ws.Cells.Item(107, 14).Value = -4638.75

# Sheet 3, Row 134
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 287182.72
$ws.Cells.Item(134, 9).Value = 1276.6
$ws.Cells.Item(134, 10).Value = 2002619.4
$ws.Cells.Item(134, 11).Value = 3829.8
$ws.Cells.Item(134, 12).Value = 6007858.199999999
$ws.Cells.Item(134, 13).Value = -1294.8
$ws.Cells.Item(134, 14).Value = -6012928.199999999

# Sheet 4, Row 16
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 1805.5
$ws.Cells.Item(16, 9).Value = 1654.5
$ws.Cells.Item(16, 10).Value = 1956.5
$ws.Cells.Item(16, 11).Value = 1654.5
$ws.Cells.Item(16, 12).Value = 1956.5
$ws.Cells.Item(16, 13).Value = -1367.5
$ws.Cells.Item(16, 14).Value = -2530.5

# Sheet 4, Row 31
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 250103.31
$ws.Cells.Item(31, 9).Value = 8128.375
$ws.Cells.Item(31, 10).Value = 482399.25
$ws.Cells.Item(31, 11).Value = 8128.375
$ws.Cells.Item(31, 12).Value = 482399.25
$ws.Cells.Item(31, 13).Value = -7833.375
$ws.Cells.Item(31, 14).Value = -482989.25

# Sheet 4, Row 34
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(34, 8).Value = 250103.31
$ws.Cells.Item(34, 9).Value = 8128.375
$ws.Cells.Item(34, 10).Value = 482399.25
$ws.Cells.Item(34, 11).Value = 8128.375
$ws.Cells.Item(34, 12).Value = 482399.25
$ws.Cells.Item(34, 13).Value = -7926.375
$ws.Cells.Item(34, 14).Value = -482803.25

# Sheet 4, Row 69
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(69, 8).Value = 49333.332
$ws.Cells.Item(69, 9).Value = 16500
$ws.Cells.Item(69, 10).Value = 115000
$ws.Cells.Item(69, 11).Value = 16500
$ws.Cells.Item(69, 12).Value = 115000
$ws.Cells.Item(69, 13).Value = -15751
$ws.Cells.Item(69, 14).Value = -116498

# Sheet 4, Row 72
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(72, 8).Value = 49333.332
$ws.Cells.Item(72, 9).Value = 16500
$ws.Cells.Item(72, 10).Value = 115000
$ws.Cells.Item(72, 11).Value = 49500
$ws.Cells.Item(72, 12).Value = 345000
$ws.Cells.Item(72, 13).Value = -45756
$ws.Cells.Item(72, 14).Value = -352488

# Sheet 4, Row 113
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(113, 8).Value = 1805.5
$ws.Cells.Item(113, 9).Value = 1654.5
$ws.Cells.Item(113, 10).Value = 1956.5
$ws.Cells.Item(113, 11).Value = 1654.5
$ws.Cells.Item(113, 12).Value = 1956.5
$ws.Cells.Item(113, 13).Value = 515.5
$ws.Cells.Item(113, 14).Value = -6296.5

# Sheet 4, Row 132
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(132, 8).Value = 2212.0967
$ws.Cells.Item(132, 9).Value = 2270.9285
$ws.Cells.Item(132, 10).Value = 1663
$ws.Cells.Item(132, 11).Value = 6812.7855
$ws.Cells.Item(132, 12).Value = 4989
$ws.Cells.Item(132, 13).Value = -4282.7855
$ws.Cells.Item(132, 14).Value = -10049

# Sheet 4, Row 134
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134, 8).Value = 3734.682
$ws.Cells.Item(134, 9).Value = 2572.875
$ws.Cells.Item(134, 11).Value = 7718.625
$ws.Cells.Item(134, 13).Value = -5183.625

# Sheet 5, Row 14
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(14, 8).Value = 1000
$ws.Cells.Item(14, 9).Value = 1000
$ws.Cells.Item(14, 11).Value = 3000

# Sheet 5, Row 68
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(68, 8).Value = 2390.9375
$ws.Cells.Item(68, 9).Value = 1971
$ws.Cells.Item(68, 11).Value = 5913
$ws.Cells.Item(68, 13).Value = -5102

# Sheet 5, Row 71
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(71, 8).Value = 2390.9375
$ws.Cells.Item(71, 9).Value = 1971
$ws.Cells.Item(71, 11).Value = 17739
$ws.Cells.Item(71, 13).Value = -13683

# Sheet 5, Row 86
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(86, 8).Value = 799.4286
$ws.Cells.Item(86, 10).Value = 782.6667
$ws.Cells.Item(86, 12).Value = 2348.0001
$ws.Cells.Item(86, 14).Value = -4720.0001

# Sheet 5, Row 89
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(89, 8).Value = 799.4286
$ws.Cells.Item(89, 10).Value = 782.6667
$ws.Cells.Item(89, 12).Value = 7044.0003
$ws.Cells.Item(89, 14).Value = -18900.0003

# Sheet 5, Row 114
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(114, 8).Value = 981
$ws.Cells.Item(114, 10).Value = 1212
$ws.Cells.Item(114, 12).Value = 3636

# Sheet 5, Row 117
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(117, 8).Value = 1985.5834
$ws.Cells.Item(117, 9).Value = 1282.25
$ws.Cells.Item(117, 11).Value = 3846.75
$ws.Cells.Item(117, 13).Value = -404.75

# Sheet 5, Row 129
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(129, 8).Value = 41672720
$ws.Cells.Item(129, 9).Value = 7570
$ws.Cells.Item(129, 10).Value = 166668180
$ws.Cells.Item(129, 11).Value = 22710
$ws.Cells.Item(129, 12).Value = 500004540
$ws.Cells.Item(129, 13).Value = -17710
$ws.Cells.Item(129, 14).Value = -500014540

# Sheet 5, Row 131
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(131, 8).Value = 4488.4927
$ws.Cells.Item(131, 10).Value = 4418
$ws.Cells.Item(131, 12).Value = 13254
$ws.Cells.Item(131, 14).Value = -23334

# Sheet 5, Row 132
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(132, 8).Value = 1960.24
$ws.Cells.Item(132, 9).Value = 1572.5714
$ws.Cells.Item(132, 11).Value = 14153.1426
$ws.Cells.Item(132, 13).Value = -11623.1426

# Sheet 5, Row 134
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(134, 8).Value = 5932.7837
$ws.Cells.Item(134, 9).Value = 1610.9231
$ws.Cells.Item(134, 11).Value = 4832.7693
$ws.Cells.Item(134, 13).Value = 237.2307000000001

# Sheet 5, Row 139
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(139, 8).Value = 1916.6383
$ws.Cells.Item(139, 9).Value = 1276.0667
$ws.Cells.Item(139, 11).Value = 3828.2001
$ws.Cells.Item(139, 13).Value = 1311.7999

# Sheet 5, Row 140
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(140, 8).Value = 117803.54
$ws.Cells.Item(140, 9).Value = 117803.54
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 353410.62
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = -348230.62

# Sheet 6, Row 2
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 267.44446
$ws.Cells.Item(2, 9).Value = 32.857143
$ws.Cells.Item(2, 10).Value = 416.72726
$ws.Cells.Item(2, 11).Value = 32.857143
$ws.Cells.Item(2, 12).Value = 416.72726
$ws.Cells.Item(2, 13).Value = 80.14285699999999
$ws.Cells.Item(2, 14).Value = -642.72726

# Sheet 6, Row 5
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(5, 8).Value = 184999.75
$ws.Cells.Item(5, 9).Value = 184999.75
$ws.Cells.Item(5, 11).Value = 184999.75
$ws.Cells.Item(5, 13).Value = -184887.75

# Sheet 6, Row 80
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(80, 8).Value = 2739.2
$ws.Cells.Item(80, 9).Value = 2352
$ws.Cells.Item(80, 11).Value = 2352
$ws.Cells.Item(80, 13).Value = -1354

# Sheet 6, Row 83
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(83, 8).Value = 2739.2
$ws.Cells.Item(83, 9).Value = 2352
$ws.Cells.Item(83, 11).Value = 11760
$ws.Cells.Item(83, 13).Value = -6768

# Sheet 6, Row 92
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(92, 8).Value = 23596.2
$ws.Cells.Item(92, 10).Value = 23596.2
$ws.Cells.Item(92, 12).Value = 23596.2
$ws.Cells.Item(92, 14).Value = -27340.2

# Sheet 6, Row 109
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(109, 8).Value = 45227.332
$ws.Cells.Item(109, 10).Value = 45227.332
$ws.Cells.Item(109, 12).Value = 45227.332
$ws.Cells.Item(109, 14).Value = -47307.332

# Sheet 6, Row 122
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(122, 8).Value = 1767.45
$ws.Cells.Item(122, 9).Value = 1334.1578
$ws.Cells.Item(122, 11).Value = 4002.4734
$ws.Cells.Item(122, 13).Value = -1552.4734

# Sheet 6, Row 126
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(126, 8).Value = 15430.375
$ws.Cells.Item(126, 9).Value = 16920.428
$ws.Cells.Item(126, 11).Value = 50761.284
$ws.Cells.Item(126, 13).Value = -48291.284

# Sheet 6, Row 132
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132, 8).Value = 22732716
$ws.Cells.Item(132, 9).Value = 29414478
$ws.Cells.Item(132, 10).Value = 14727.6
$ws.Cells.Item(132, 11).Value = 88243434
$ws.Cells.Item(132, 12).Value = 44182.8
$ws.Cells.Item(132, 13).Value = -88240904
$ws.Cells.Item(132, 14).Value = -49242.8

# Sheet 6, Row 134
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(134, 8).Value = 99989
$ws.Cells.Item(134, 10).Value = 99989
$ws.Cells.Item(134, 12).Value = 299967
$ws.Cells.Item(134, 14).Value = -305037

# Sheet 7, Row 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 106900.7
$ws.Cells.Item(7, 9).Value = 5624.5
$ws.Cells.Item(7, 11).Value = 5624.5
$ws.Cells.Item(7, 13).Value = -5512.5

# Sheet 7, Row 36
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(36, 8).Value = 78180
$ws.Cells.Item(36, 10).Value = 78180
$ws.Cells.Item(36, 12).Value = 78180
$ws.Cells.Item(36, 14).Value = -79304

# Sheet 7, Row 40
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 3346.889
$ws.Cells.Item(40, 9).Value = 2557.2307
$ws.Cells.Item(40, 11).Value = 2557.2307
$ws.Cells.Item(40, 13).Value = -2421.2307

# Sheet 7, Row 42
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(42, 8).Value = 16514
$ws.Cells.Item(42, 10).Value = 8028
$ws.Cells.Item(42, 12).Value = 8028

# Sheet 7, Row 46
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(46, 8).Value = 14850.919
$ws.Cells.Item(46, 9).Value = 2125.25
$ws.Cells.Item(46, 11).Value = 2125.25
$ws.Cells.Item(46, 13).Value = -1937.25

# Sheet 7, Row 49
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(49, 8).Value = 16514
$ws.Cells.Item(49, 10).Value = 8028
$ws.Cells.Item(49, 12).Value = 8028

# Sheet 7, Row 55
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(55, 8).Value = 29412124
$ws.Cells.Item(55, 10).Value = 440.27274
$ws.Cells.Item(55, 12).Value = 440.27274
$ws.Cells.Item(55, 14).Value = -786.27274

# Sheet 7, Row 68
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(68, 8).Value = 3436.6924
$ws.Cells.Item(68, 9).Value = 3096.3333
$ws.Cells.Item(68, 10).Value = 3728.4285
$ws.Cells.Item(68, 11).Value = 3096.3333
$ws.Cells.Item(68, 12).Value = 3728.4285
$ws.Cells.Item(68, 13).Value = -2347.3333
$ws.Cells.Item(68, 14).Value = -5226.4285

# Sheet 7, Row 71
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(71, 8).Value = 3436.6924
$ws.Cells.Item(71, 9).Value = 3096.3333
$ws.Cells.Item(71, 10).Value = 3728.4285
$ws.Cells.Item(71, 11).Value = 15481.6665
$ws.Cells.Item(71, 12).Value = 18642.1425
$ws.Cells.Item(71, 13).Value = -11737.6665
$ws.Cells.Item(71, 14).Value = -26130.1425

# Sheet 7, Row 82
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(82, 8).Value = 3697.2222
$ws.Cells.Item(82, 9).Value = 4045.8333
$ws.Cells.Item(82, 11).Value = 4045.8333
$ws.Cells.Item(82, 13).Value = -3684.8333

# Sheet 7, Row 85
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(85, 8).Value = 3697.2222
$ws.Cells.Item(85, 9).Value = 4045.8333
$ws.Cells.Item(85, 11).Value = 4045.8333
$ws.Cells.Item(85, 13).Value = -2797.8333

# Sheet 7, Row 124
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(124, 8).Value = 109995
$ws.Cells.Item(124, 10).Value = 109995
$ws.Cells.Item(124, 12).Value = 109995
$ws.Cells.Item(124, 14).Value = -119815

# Sheet 7, Row 126
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(126, 8).Value = 106900.7
$ws.Cells.Item(126, 9).Value = 5624.5
$ws.Cells.Item(126, 11).Value = 16873.5
$ws.Cells.Item(126, 13).Value = -14403.5

# Sheet 7, Row 132
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(132, 8).Value = 203890.08
$ws.Cells.Item(132, 9).Value = 4257.8
$ws.Cells.Item(132, 11).Value = 12773.4
$ws.Cells.Item(132, 13).Value = -10243.4

# Sheet 7, Row 136
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(136, 8).Value = 50108.04
$ws.Cells.Item(136, 9).Value = 6640.7
$ws.Cells.Item(136, 11).Value = 19922.1
$ws.Cells.Item(136, 13).Value = -17372.1

# Sheet 8, Row 46
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(46, 8).Value = 82799.6
$ws.Cells.Item(46, 10).Value = 82799.6
$ws.Cells.Item(46, 12).Value = 82799.6
$ws.Cells.Item(46, 14).Value = -83261.6

# Sheet 8, Row 62
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(62, 8).Value = 16672283
$ws.Cells.Item(62, 9).Value = 4300.6665
$ws.Cells.Item(62, 11).Value = 4300.6665
$ws.Cells.Item(62, 13).Value = -3676.6665

# Sheet 8, Row 65
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(65, 8).Value = 16672283
$ws.Cells.Item(65, 9).Value = 4300.6665
$ws.Cells.Item(65, 11).Value = 21503.3325
$ws.Cells.Item(65, 13).Value = -18383.3325

# Sheet 8, Row 96
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(96, 8).Value = 5399.727
$ws.Cells.Item(96, 10).Value = 5199.4
$ws.Cells.Item(96, 12).Value = 5199.4
$ws.Cells.Item(96, 14).Value = -7945.4

# Sheet 8, Row 105
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(105, 8).Value = 22521.8
$ws.Cells.Item(105, 10).Value = 22521.8
$ws.Cells.Item(105, 12).Value = 22521.8
$ws.Cells.Item(105, 14).Value = -29509.8

# Sheet 8, Row 107
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(107, 8).Value = 20834532
$ws.Cells.Item(107, 9).Value = 31251468
$ws.Cells.Item(107, 11).Value = 93754404
$ws.Cells.Item(107, 13).Value = -93752484

# Sheet 8, Row 122
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(122, 8).Value = 260
$ws.Cells.Item(122, 9).Value = 260
$ws.Cells.Item(122, 11).Value = 780

# Sheet 8, Row 132
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 2535.5535
$ws.Cells.Item(132, 9).Value = 1957.0817
$ws.Cells.Item(132, 11).Value = 5871.2451
$ws.Cells.Item(132, 13).Value = -3341.2451

# Sheet 8, Row 134
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(134, 8).Value = 82799.6
$ws.Cells.Item(134, 10).Value = 82799.6
$ws.Cells.Item(134, 12).Value = 248398.8
$ws.Cells.Item(134, 14).Value = -253468.8

# Sheet 8, Row 135
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(135, 8).Value = 126333
$ws.Cells.Item(135, 10).Value = 126333
$ws.Cells.Item(135, 12).Value = 126333
$ws.Cells.Item(135, 14).Value = -136473

# Sheet 8, Row 136
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(136, 8).Value = 2465.0889
$ws.Cells.Item(136, 9).Value = 2477.9697
$ws.Cells.Item(136, 10).Value = 2429.6667
$ws.Cells.Item(136, 11).Value = 7433.909100000001
$ws.Cells.Item(136, 12).Value = 7289.000100000001
$ws.Cells.Item(136, 13).Value = -4883.909100000001
$ws.Cells.Item(136, 14).Value = -12389.0001

# Special additions (new cells)
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(14, 13).Value = -2827
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(114, 14).Value = -10144
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(42, 14).Value = -9154
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(49, 14).Value = -8322
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(122, 13).Value = 1670

# Special removals (clear cell entirely)
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(140, 14).ClearContents()
